$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B3").Value = "2.1.0"
$wsMeta.Range("B8").Value = "2025-12-19T08:44:55+00:00"

$wsVS = $wb.Worksheets.Item("Include ValueSet #0")
$wsVS.Range("A2").Value = "https://mos.esante.gouv.fr/NOS/JDV_J283-PrestationsIndirects_SERAFIN/FHIR/JDV-J283-PrestationsIndirects-SERAFIN"
$wsVS.Range("A3").Value = "https://mos.esante.gouv.fr/NOS/JDV_J284-PrestationsDirects_SERAFIN/FHIR/JDV-J284-PrestationsDirects-SERAFIN"
